$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = $ws.Range("A11").Value2
$ws.Range("B12").Value = $ws.Range("B11").Value2
$ws.Range("C12").Value = $ws.Range("C11").Value2
$ws.Range("D12").Value = $ws.Range("D11").Value2
